$wb = $excel.ActiveWorkbook

# hunk 0: ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1145.1111
$ws.Range("I92").Value = 642.44446
$ws.Range("J92").Value = 2150.4443
$ws.Range("K92").Value = 642.44446
$ws.Range("L92").Value = 2150.4443
$ws.Range("M92").Value = 605.55554
$ws.Range("N92").Value = -4646.4443

# hunk 1: ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6669170
$ws.Range("I137").Value = 12501725
$ws.Range("J137").Value = 3392.8572
$ws.Range("K137").Value = 37505175
$ws.Range("L137").Value = 10178.5716
$ws.Range("M137").Value = -37502625
$ws.Range("N137").Value = -15278.5716

# hunk 2: ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16057.494
$ws.Range("I32").Value = 17312.016
$ws.Range("K32").Value = 17312.016
$ws.Range("M32").Value = -17025.016

# hunk 3: ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34553224
$ws.Range("I61").Value = 43523080
$ws.Range("J61").Value = 168788.83
$ws.Range("K61").Value = 43523080
$ws.Range("L61").Value = 168788.83
$ws.Range("M61").Value = -43522868
$ws.Range("N61").Value = -169212.83

# hunk 4: ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5041049.5
$ws.Range("I74").Value = 6973164
$ws.Range("J74").Value = 72756.42999999999
$ws.Range("K74").Value = 6973164
$ws.Range("L74").Value = 72756.42999999999
$ws.Range("M74").Value = -6972290
$ws.Range("N74").Value = -74504.42999999999

# hunk 5: ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5041049.5
$ws.Range("I77").Value = 6973164
$ws.Range("J77").Value = 72756.42999999999
$ws.Range("K77").Value = 34865820
$ws.Range("L77").Value = 363782.15
$ws.Range("M77").Value = -34861452
$ws.Range("N77").Value = -372518.15

# hunk 6: ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 69498.03
$ws.Range("I132").Value = 53896.895
$ws.Range("J132").Value = 94199.836
$ws.Range("K132").Value = 161690.685
$ws.Range("L132").Value = 282599.508
$ws.Range("M132").Value = -159160.685
$ws.Range("N132").Value = -287659.508

# hunk 7: ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 34553224
$ws.Range("I136").Value = 43523080
$ws.Range("J136").Value = 168788.83
$ws.Range("K136").Value = 130569240
$ws.Range("L136").Value = 506366.49
$ws.Range("M136").Value = -130566690
$ws.Range("N136").Value = -511466.49

# hunk 8: BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1298.75
$ws.Range("I94").Value = 722.75
$ws.Range("J94").Value = 1586.75
$ws.Range("K94").Value = 722.75
$ws.Range("L94").Value = 1586.75
$ws.Range("M94").Value = -271.75
$ws.Range("N94").Value = -2488.75

# hunk 9: BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1211.4
$ws.Range("I134").Value = 1155.0714
$ws.Range("K134").Value = 3465.2142
$ws.Range("M134").Value = -930.2142000000003

# hunk 10: CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3152.3103
$ws.Range("I31").Value = 1521.6
$ws.Range("J31").Value = 4899.5
$ws.Range("K31").Value = 1521.6
$ws.Range("L31").Value = 4899.5
$ws.Range("M31").Value = -1226.6
$ws.Range("N31").Value = -5489.5

# hunk 11: CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3152.3103
$ws.Range("I34").Value = 1521.6
$ws.Range("J34").Value = 4899.5
$ws.Range("K34").Value = 1521.6
$ws.Range("L34").Value = 4899.5
$ws.Range("M34").Value = -1319.6
$ws.Range("N34").Value = -5303.5

# hunk 12: CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 24391994
$ws.Range("I58").Value = 34483800
$ws.Range("J58").Value = 3460.25
$ws.Range("K58").Value = 34483800
$ws.Range("L58").Value = 3460.25
$ws.Range("M58").Value = -34483597
$ws.Range("N58").Value = -3866.25

# hunk 13: CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 37914.895
$ws.Range("I132").Value = 2033.2222
$ws.Range("J132").Value = 102501.9
$ws.Range("K132").Value = 6099.6666
$ws.Range("L132").Value = 307505.7
$ws.Range("M132").Value = -3569.6666
$ws.Range("N132").Value = -312565.7

# hunk 14: CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 47501.047
$ws.Range("I134").Value = 2247.5881
$ws.Range("J134").Value = 201362.8
$ws.Range("K134").Value = 6742.7643
$ws.Range("L134").Value = 604088.3999999999
$ws.Range("M134").Value = -4207.7643
$ws.Range("N134").Value = -609158.3999999999

# hunk 15: CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 24391994
$ws.Range("I136").Value = 34483800
$ws.Range("J136").Value = 3460.25
$ws.Range("K136").Value = 103451400
$ws.Range("L136").Value = 10380.75
$ws.Range("M136").Value = -103448850
$ws.Range("N136").Value = -15480.75

# hunk 16: CUL row 110
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 5664.25
$ws.Range("I110").Value = 1363.5
$ws.Range("J110").Value = 9965
$ws.Range("K110").Value = 4090.5
$ws.Range("L110").Value = 29895
$ws.Range("M110").Value = -0.5
$ws.Range("N110").Value = -38075

# hunk 17: CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 831.38464
$ws.Range("I122").Value = 552.3333
$ws.Range("J122").Value = 1070.5714
$ws.Range("K122").Value = 4970.9997
$ws.Range("L122").Value = 9635.142600000001
$ws.Range("M122").Value = -2520.9997
$ws.Range("N122").Value = -14535.1426

# hunk 18: CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1044
$ws.Range("I132").Value = 653.43475
$ws.Range("K132").Value = 5880.91275
$ws.Range("M132").Value = -3350.91275

# hunk 19: CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2266.2124
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 2266.2124
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 6798.637199999999
$ws.Range("N140").Value = -17158.6372
$ws.Range("M140").ClearContents()

# hunk 20: GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 82853.16
$ws.Range("I132").Value = 73480.5
$ws.Range("K132").Value = 220441.5
$ws.Range("M132").Value = -217911.5

# hunk 21: LTW row 44
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 5000
$ws.Range("L44").Value = 5000
$ws.Range("N44").Value = -5912

# hunk 22: LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2200
$ws.Range("I61").Value = 2200
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2200
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1998
$ws.Range("N61").ClearContents()

# hunk 23: LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1304.875
$ws.Range("I93").Value = 920.1
$ws.Range("J93").Value = 1579.7142
$ws.Range("K93").Value = 920.1
$ws.Range("L93").Value = 1579.7142
$ws.Range("M93").Value = 327.9
$ws.Range("N93").Value = -4075.7142

# hunk 24: LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2200
$ws.Range("I113").Value = 2200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2200
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -30
$ws.Range("N113").ClearContents()

# hunk 25: LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 81127.766
$ws.Range("I132").Value = 2942
$ws.Range("J132").Value = 95343.37
$ws.Range("K132").Value = 8826
$ws.Range("L132").Value = 286030.11
$ws.Range("M132").Value = -6296
$ws.Range("N132").Value = -291090.11

# hunk 26: LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 334314.16
$ws.Range("I136").Value = 500990
$ws.Range("J136").Value = 250976.25
$ws.Range("K136").Value = 1502970
$ws.Range("L136").Value = 752928.75
$ws.Range("M136").Value = -1500420
$ws.Range("N136").Value = -758028.75

# hunk 27: WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 92436.23
$ws.Range("I132").Value = 56588.832
$ws.Range("K132").Value = 169766.496
$ws.Range("M132").Value = -167236.496

# hunk 28: WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 106353.266
$ws.Range("I136").Value = 72400.78999999999
$ws.Range("K136").Value = 217202.37
$ws.Range("M136").Value = -214652.37
